$d = $word.ActiveDocument

# --- Locate the paragraph that holds "offline feladatok ... ID101 ..."
# (the new bullet must be inserted right after it) -----------------------

$offlinePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt -like "*offline feladatok*ID101*") {
        $offlinePara = $d.Paragraphs($i)
    }
}
if ($offlinePara -eq $null) {
    throw "Could not find the 'offline feladatok ... ID101 ...' paragraph"
}

# --- Remove the _GoBack bookmark that currently sits after "(Bálint)" ---

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Insert a brand-new list paragraph right after the "offline
# feladatok ..." bullet. It gets the same paragraph style / numbering
# (InsertParagraphAfter clones the source paragraph's pPr/rPr), then its
# content is populated precisely via InsertXML so the spell-check markers
# around "életkorilag" and the relocated _GoBack bookmark land exactly
# where the target revision puts them (right after the last run, before
# the paragraph end). ------------------------------------------------

$offlinePara.Range.InsertParagraphAfter()
$newPara = $offlinePara.Next()

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr>' +
    '<w:pStyle w:val="Listaszerbekezds"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
    '</w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr>' +
    '<w:t>életkorilag</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr>' +
    '<w:t xml:space="preserve"> nagyon heterogén, nincsenek fiatal felnőttek</w:t>' +
    '</w:r>' +
    '<w:bookmarkStart w:id="6" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="6"/>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$newPara.Range.InsertXML($xml)

Write-Output "done"
